# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Column F ("tamano-centro") metadata block (rows 2-4) and the
# now-obsolete mapping-file reference in row 5 are refreshed to point at
# the newly curated iaest-measure dimension; likewise column I
# ("direccion-provincial-nombre").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F: tamano-centro ---------------------------------------
$ws.Range("F2").Value = "iaest-measure:tamano-centro"
$ws.Range("F3").Value = "medida"
$ws.Range("F4").Value = "xsd:int"
$ws.Range("F5").Clear()

# --- Column I: direccion-provincial-nombre --------------------------
$ws.Range("I2").Value = "iaest-measure:direccion-provincial-nombre"
$ws.Range("I3").Value = "medida"
$ws.Range("I4").Value = "xsd:int"
